# Update the "Metadata" sheet of the ValueSet workbook to the new release:
#  - Version 5.0.0 -> 6.0.0
#  - Date 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
#  - Publisher value filled in ("Alvearie Team")
#  - Duplicate "Contact" / "No display for ContactDetail" rows replaced with a
#    single "Jurisdiction" / "United States of America" row (net: one row removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version (row 3)
$ws.Range("B3").Value = "6.0.0"

# Date (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (row 9) - was empty
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 currently holds "Contact" / "No display for ContactDetail" and row 11
# is an exact duplicate of it. Delete row 11 (the duplicate) and then turn the
# remaining row 10 into the new "Jurisdiction" / "United States of America" row.
$ws.Rows.Item(11).Delete()

$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
